# React Redux with comments
# Cleans up stray w:proofErr spell/grammar-check markers (and the run
# splits they forced) across several list items, and appends a new
# "CSS or flex box course UDEMY" line at the end of the document.

$d = $word.ActiveDocument

function Set-CleanParagraph($Index, $Text) {
    $para = $d.Paragraphs($Index)
    $escaped = $Text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
           '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body>' +
           '<w:p>' +
           '<w:pPr><w:rPr><w:lang w:val="en-IN"/></w:rPr></w:pPr>' +
           '<w:r><w:rPr><w:lang w:val="en-IN"/></w:rPr><w:t>' + $escaped + '</w:t></w:r>' +
           '</w:p>' +
           '</w:body></w:document>' +
           '</pkg:xmlData></pkg:part></pkg:package>'
    $para.Range.InsertXML($xml)
}

# Paragraph 1: "Iteraotr & Generators." - merge runs, drop proofErr
Set-CleanParagraph 1 "Iteraotr & Generators."

# Paragraph 5: "Generators(Ffunction*)" - merge runs, drop proofErr
Set-CleanParagraph 5 "Generators(Ffunction*)"

# Paragraph 6: "Bind()" - drop proofErr (text already single run)
Set-CleanParagraph 6 "Bind()"

# Paragraph 7: "Call()" - drop proofErr
Set-CleanParagraph 7 "Call()"

# Paragraph 8: "Apply()" - drop proofErr
Set-CleanParagraph 8 "Apply()"

# Paragraph 10: "Concat()" - merge runs, drop proofErr
Set-CleanParagraph 10 "Concat()"

# Paragraph 11: "Match()" - drop proofErr
Set-CleanParagraph 11 "Match()"

# Append a new paragraph after the trailing empty paragraph with the
# new list item text.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs($d.Paragraphs.Count)
$newPara.Range.Text = "CSS or flex box course UDEMY"
